$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the name values in column B (shared strings: mert->Mert, Hasan->Ali, AHMET->Ayse) ---
$ws.Range("B2").Value = "Mert"
$ws.Range("B3").Value = "Ali"
$ws.Range("B4").Value = "Ayşe"

# --- Update the ID values in column A for rows 3 and 4 ---
$ws.Range("A3").Value = 5678
$ws.Range("A4").Value = 9012

# --- Apply bold + centered formatting to the header row (A1:C1) ---
# Build the combined style on a single cell first (avoids a transient
# "bold only" style from lingering in the style table), then propagate
# it to the rest of the header via a format-only copy/paste.
$header = $ws.Cells.Item(1, 1)
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108

$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
